# Apply cryptocurrency price/volume updates per commit
# "Updated cryptos list on Sun Oct  8 06:54:05 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as numbers by Excel;
# force Text format first so the literal string (incl. trailing zeros) is kept.
$textForceCells = @("D5", "D8", "D9", "D11", "D14", "D15", "D16", "D18", "D20", "D22", "D23", "D24", "D25", "D30", "D35", "D40", "D44", "D48", "D50", "D51")
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "27.918.25"
$ws.Range("E2").Value = "  -0.01%  "

# Row 3
$ws.Range("D3").Value = "1.630.16"
$ws.Range("E3").Value = "  -0.60%  "

# Row 4
$ws.Range("E4").Value = "  +0.05%  "

# Row 5
$ws.Range("D5").Value = "211.80"

# Row 6
$ws.Range("E6").Value = "  -0.32%  "

# Row 7
$ws.Range("E7").Value = "  +0.05%  "

# Row 8
$ws.Range("D8").Value = "23.35"
$ws.Range("E8").Value = "  -1.20%  "

# Row 9
$ws.Range("D9").Value = "0.257"
$ws.Range("E9").Value = "  -1.96%  "

# Row 10
$ws.Range("E10").Value = "  -0.27%  "

# Row 11
$ws.Range("D11").Value = "0.0880"
$ws.Range("E11").Value = "  +0.63%  "

# Row 12
$ws.Range("D12").Value = "1.862.75"
$ws.Range("E12").Value = "  -0.53%  "

# Row 13
$ws.Range("D13").Value = "1.626.91"
$ws.Range("E13").Value = "  -0.80%  "

# Row 14
$ws.Range("D14").Value = "4.04"
$ws.Range("E14").Value = "  -1.27%  "

# Row 15
$ws.Range("D15").Value = "0.561"
$ws.Range("E15").Value = "  -2.37%  "

# Row 16
$ws.Range("D16").Value = "65.56"
$ws.Range("E16").Value = "  -1.01%  "

# Row 17
$ws.Range("D17").Value = "27.927.98"
$ws.Range("E17").Value = "  +0.07%  "

# Row 18
$ws.Range("D18").Value = "230.46"
$ws.Range("E18").Value = "  -0.87%  "

# Row 19
$ws.Range("E19").Value = "  -0.02%  "

# Row 20
$ws.Range("D20").Value = "7.64"
$ws.Range("E20").Value = "  +0.69%  "

# Row 21
$ws.Range("E21").Value = "  -0.07%  "

# Row 22
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "4.34"
$ws.Range("E22").Value = "  -0.66%  "

# Row 23
$ws.Range("B23").Value = "Avalanche"
$ws.Range("C23").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D23").Value = "10.29"
$ws.Range("E23").Value = "  -4.95%  "

# Row 24
$ws.Range("D24").Value = "2.02"
$ws.Range("E24").Value = "  -1.80%  "

# Row 25
$ws.Range("D25").Value = "154.89"
$ws.Range("E25").Value = "  +2.01%  "

# Row 27
$ws.Range("E27").Value = "  -0.11%  "

# Row 28
$ws.Range("E28").Value = "  -1.15%  "

# Row 29
$ws.Range("E29").Value = "  +0.02%  "

# Row 30
$ws.Range("D30").Value = "1.18"
$ws.Range("E30").Value = "  -0.51%  "

# Row 31
$ws.Range("E31").Value = "  -0.52%  "

# Row 32
$ws.Range("E32").Value = "  +1.76%  "

# Row 33
$ws.Range("D33").Value = "1.400.57"
$ws.Range("E33").Value = "  -1.08%  "

# Row 34
$ws.Range("E34").Value = "  -1.40%  "

# Row 35
$ws.Range("D35").Value = "1.57"
$ws.Range("E35").Value = "  +0.09%  "

# Row 36
$ws.Range("E36").Value = "  +10.97%  "

# Row 37
$ws.Range("E37").Value = "  +0.39%  "

# Row 38
$ws.Range("E38").Value = "  +2.06%  "

# Row 39
$ws.Range("E39").Value = "  +0.11%  "

# Row 40
$ws.Range("D40").Value = "0.864"
$ws.Range("E40").Value = "  -3.13%  "

# Row 41
$ws.Range("E41").Value = "  -0.25%  "

# Row 42
$ws.Range("E42").Value = "  -0.02%  "

# Row 43
$ws.Range("E43").Value = "  +1.03%  "

# Row 44
$ws.Range("D44").Value = "66.43"
$ws.Range("E44").Value = "  -0.63%  "

# Row 45
$ws.Range("E45").Value = "  +0.33%  "

# Row 46
$ws.Range("E46").Value = "  -0.15%  "

# Row 47
$ws.Range("D47").Value = "1.773.00"
$ws.Range("E47").Value = "  -0.43%  "

# Row 48
$ws.Range("D48").Value = "88.23"
$ws.Range("E48").Value = "  -0.06%  "

# Row 49
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.0₆0103"
$ws.Range("E49").Value = "  -1.91%  "

# Row 50
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "0.100"
$ws.Range("E50").Value = "  -0.32%  "

# Row 51
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "0.0505"
$ws.Range("E51").Value = "  -0.16%  "

